# Update "Run 0" (column B) results for Station 2 best results sheet
# Source: table_model____spi-3__lasso_($\beta_=_$2.0).xlsx re-run values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$runZeroValues = @{
    2 = 0.8833950785914961
    3 = 0.7736611616623972
    4 = 0.2105518460493915
    5 = -0.3221709317925977
    6 = 0.2149251723985612
    7 = -0.5556526060951155
    8 = -0.2622258882598426
    9 = 0.4594347466302878
    10 = 0.9007745340493207
    11 = 0.5377056971462364
    12 = -0.865507539493848
    13 = -0.5632622355798349
    14 = -0.5623665235252416
    15 = -0.4981542579122668
    16 = 0.9423025476386531
    17 = 1.186902581930848
    18 = 1.124028253452189
    19 = -1.173361664834818
    20 = -1.303792812116957
    21 = -1.257076136507502
    22 = -0.3979610366483977
    23 = -1.325454284204314
    24 = -1.663241938363834
    25 = -1.160389593668054
    26 = -1.140330368997835
    27 = -0.2019427402850872
    28 = -1.256466064811198
    29 = -0.6775301855129591
    30 = 0.3384905666306841
    31 = 1.124474408949146
    32 = 0.8548992116866377
    33 = -0.3276873662770162
    34 = -0.725223773105192
    35 = -0.493428640959891
    36 = -0.2430999785521488
    37 = -0.8101624106440117
    38 = -0.9557642662274581
    39 = -0.9670371426900926
    40 = 0.3758301326153209
    41 = 0.3835135775911769
    42 = 0.7225982252851023
    43 = -0.2698414022581147
    44 = 0.5137086292164884
    45 = 1.594028089333894
    46 = 2.08347934806333
    47 = 2.344772587040754
    48 = 0.8393974088420543
    49 = -0.1256096277171001
    50 = -0.955425990945407
    51 = -0.8604335994521285
    52 = -0.7216021831500352
    53 = -1.002862513697061
    54 = -0.2564001756977966
    55 = 0.8733550774898844
    56 = 1.618383604650495
    57 = 1.827627236028162
    58 = 1.34529538913001
    59 = 0.7289340008002547
    60 = -0.4940084644044236
    61 = 0.175366130815175
    62 = 0.688201845923031
    63 = 0.8984675767298729
    64 = -0.127763737510738
    65 = 1.869754431625355
    66 = 1.949028540012478
    67 = 2.155690906379023
    68 = 0.4701190424666771
    69 = 0.3903288936073673
    70 = 0.01276558797508542
    71 = 0.09420807993398615
    72 = 0.7152192026080975
    73 = 0.6230233909843215
    74 = 0.6463247136355199
    75 = 0.3384722415489917
    76 = -0.1719918537984761
    77 = 0.4370149967340428
    78 = -0.4595275989506945
    79 = -0.01988028872083243
    80 = 0.6248252517479855
    81 = 1.328724940832351
    82 = 1.234484010956354
    83 = -0.3948208933108187
    84 = -0.2049195052343094
    85 = -0.9148164282347584
    86 = -0.3825898103049017
    87 = -1.489439007415337
    88 = -1.278726656137989
    89 = -0.7054162657870815
    90 = -0.593032097697648
    91 = 0.3625917913938151
    92 = 0.667245399029806
    93 = 0.697065603110913
    94 = 0.4423555056875265
    95 = 0.2737266415475333
    96 = -0.3575617979965108
    97 = -1.123778061959173
    98 = -0.7626364663500418
    99 = -0.1766820870840194
    100 = 0.05875504377207696
    101 = -0.1893661712397061
    102 = -0.4589032049653482
    103 = -1.046907065762377
    104 = -1.615934318215019
    105 = -1.841984913935895
    106 = -0.4426286129516878
    107 = -0.5292662414101115
    108 = 0.7230789784347936
    109 = 1.342715688942997
    110 = 1.604571030033696
    111 = 1.100424581959709
    112 = 1.600933974533124
    113 = 1.804221614804293
    114 = 1.695787125639826
    115 = 0.5479980053863395
    116 = 0.4520356295363708
    117 = 1.17697108482315
    118 = 1.55516476307115
    119 = 1.353133784813052
    120 = 0.5300983206323587
    121 = 1.467119297739529
    122 = 1.539167779183578
    123 = 1.922219884065255
    124 = 0.4725728045386651
    125 = 0.9254103838204238
    126 = 0.008236104198193567
    127 = -0.9106811909267245
    128 = 0.1461908068160387
    129 = 0.5873919487366224
    130 = 1.818525040523871
    131 = 0.8312934070473759
    132 = 1.172258414965312
    133 = 0.1229385674984841
    134 = 0.1179833422271376
    135 = 0.2619528087184982
    136 = 1.132321445771324
    137 = 0.6348431870438699
    138 = -0.2040123955994738
    139 = -0.4859110370450214
}

foreach ($row in $runZeroValues.Keys) {
    $ws.Range("B$row").Value = $runZeroValues[$row]
}
